# Add merged cells for fixtures.
#
# The header row currently has one cell per fixture (F1..N1). Each of those
# fixture headers needs to become a 3-column-wide merged block (keeping the
# same two blank columns as spacer columns to its right), so two new blank
# columns are inserted immediately before every fixture header cell from
# F1 through N1. Column E1 ("LEE - BOU") is left where it is and becomes the
# left edge of the first merged block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before each original fixture-header column.
# Walk right-to-left (N down to F) so that column letters used below always
# refer to their original, not-yet-shifted position.
$columnPairsToInsert = @("N:O", "M:N", "L:M", "K:L", "J:K", "I:J", "H:I", "G:H", "F:G")
foreach ($pair in $columnPairsToInsert) {
    $ws.Columns($pair).Insert()
}

# After the inserts above, the fixture headers (originally F1..N1) now sit
# three columns apart starting at H1, with E1 untouched:
#   E1, H1, K1, N1, Q1, T1, W1, Z1, AC1, AF1
# Merge each fixture header together with the two blank columns that follow
# it.
$ws.Range("E1:G1").Merge()
$ws.Range("H1:J1").Merge()
$ws.Range("K1:M1").Merge()
$ws.Range("N1:P1").Merge()
$ws.Range("Q1:S1").Merge()
$ws.Range("T1:V1").Merge()
$ws.Range("W1:Y1").Merge()
$ws.Range("Z1:AB1").Merge()
$ws.Range("AC1:AE1").Merge()
$ws.Range("AF1:AH1").Merge()

# Touch every blank header cell with a (no-op) formatting change so the used
# range / sheet dimension grows to cover the new columns up to AH1, matching
# the widened header row.
$ws.Range("F1:AH1").Font.Bold = $false

Write-Host "Merged fixture header cells E1:G1 through AF1:AH1"
